# Fruta / hortaliza, semanal
# Inserts two new weekly price records (rows 21-22) for
# "Vega Monumental Concepción" - Mango, pushing the existing
# historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right after the current row 20, shifting all
# subsequent rows (old 21..58) down to (23..60).
$ws.Range("A21:A22").EntireRow.Insert()

# New row 21
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value = 44427
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100108
$ws.Range("H21").Value = "Tropicales y subtropicales"
$ws.Range("I21").Value = 100108002
$ws.Range("J21").Value = "Mango"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 8500
$ws.Range("O21").Value = 9000
$ws.Range("P21").Value = 8750
$ws.Range("Q21").Value = "$/bandeja 4 kilos"
$ws.Range("R21").Value = "Brasil"
$ws.Range("S21").Value = 2188
$ws.Range("T21").Value = 4

# New row 22
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = "Vega Monumental Concepción"
$ws.Range("C22").Value = "Bíobío"
$ws.Range("D22").Value = 44421
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100108
$ws.Range("H22").Value = "Tropicales y subtropicales"
$ws.Range("I22").Value = 100108002
$ws.Range("J22").Value = "Mango"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 8500
$ws.Range("O22").Value = 9000
$ws.Range("P22").Value = 8750
$ws.Range("Q22").Value = "$/bandeja 4 kilos"
$ws.Range("R22").Value = "Brasil"
$ws.Range("S22").Value = 2188
$ws.Range("T22").Value = 4
